$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.469.27'
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").Value = '1.997.56'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.38'
$ws.Range("E5").Value = '  -2.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.631'
$ws.Range("E6").Value = '  -1.95%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '61.79'
$ws.Range("E7").Value = '  -1.23%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.385'
$ws.Range("E9").Value = '  +3.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '57.09'
$ws.Range("E10").Value = '  -3.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0777'
$ws.Range("E11").Value = '  +4.27%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.103'
$ws.Range("E12").Value = '  -0.94%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.12'
$ws.Range("E13").Value = '  +13.37%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.878'
$ws.Range("E14").Value = '  -3.68%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.16'
$ws.Range("E15").Value = '  -4.76%  '
$ws.Range("D16").Value = '2.290.81'
$ws.Range("E16").Value = '  -0.89%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.50'
$ws.Range("E17").Value = '  +1.10%  '
$ws.Range("D18").Value = '2.002.80'
$ws.Range("E18").Value = '  +0.53%  '
$ws.Range("D19").Value = '36.367.87'
$ws.Range("E19").Value = '  +0.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.70'
$ws.Range("E20").Value = '  -0.59%  '
$ws.Range("D21").Value = '0.0₃0871'
$ws.Range("E21").Value = '  +1.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.28'
$ws.Range("E22").Value = '  -0.62%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '234.21'
$ws.Range("E23").Value = '  -0.14%  '
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("E25").Value = '  -5.84%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.32'
$ws.Range("E26").Value = '  +0.27%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.70'
$ws.Range("E27").Value = '  +0.80%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.139'
$ws.Range("E28").Value = '  +20.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '159.49'
$ws.Range("E29").Value = '  -2.37%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.08'
$ws.Range("E30").Value = '  +2.02%  '
$ws.Range("E31").Value = '  -0.46%  '
$ws.Range("E32").Value = '  -0.89%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.95'
$ws.Range("E33").Value = '  -3.46%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0616'
$ws.Range("E34").Value = '  +0.99%  '
$ws.Range("E35").Value = '  -2.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.55'
$ws.Range("E36").Value = '  +10.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.37'
$ws.Range("E37").Value = '  -4.43%  '
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.83'
$ws.Range("E39").Value = '  +0.88%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.21'
$ws.Range("E40").Value = '  +21.61%  '
$ws.Range("B41").Value = 'Cronos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0997'
$ws.Range("E41").Value = '  -3.72%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.25'
$ws.Range("E42").Value = '  +2.42%  '
$ws.Range("E43").Value = '  -0.24%  '
$ws.Range("E44").Value = '  -1.32%  '
$ws.Range("E45").Value = '  -1.43%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.66'
$ws.Range("E46").Value = '  +0.08%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '92.49'
$ws.Range("E47").Value = '  -2.26%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.57'
$ws.Range("E48").Value = '  -3.88%  '
$ws.Range("D49").Value = '1.355.21'
$ws.Range("E49").Value = '  -5.86%  '
$ws.Range("E50").Value = '  -2.16%  '
$ws.Range("D51").Value = '2.184.39'
$ws.Range("E51").Value = '  -0.70%  '